$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.518.24"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "1.571.52"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -1.58%  "
$ws.Range("D5").Value = "211.21"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "0.0872"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "1.795.20"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "1.547.17"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "27.476.66"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "62.79"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "226.49"
$ws.Range("E18").Value = "  +5.15%  "
$ws.Range("D19").Value = "7.51"
$ws.Range("D20").Value = "0.0₃0706"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "0.993"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "9.41"
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "150.14"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").Value = "15.17"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "1.456.44"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.66"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.992"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("E44").Value = "  +7.02%  "
$ws.Range("E45").Value = "  -3.18%  "
$ws.Range("D46").Value = "63.89"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "1.706.87"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "86.93"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("D50").Value = "0.0524"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  -1.40%  "
